{"js": "// Replace the paragraph that currently reads:\n//   \"Acerca da quantidade de horas a serem alocadas, por se tratar de um\n//    projeto cargaHoraria para PRNCoordenador coordenadorPRNTxtCoordenador,\n//    paragrafo8. \"\n// with just:\n//   \"paragrafo8. \"\n// keeping the paragraph's own formatting (indent / run size, etc.) intact.\n\nconst body = context.document.body;\n\n// Locate the paragraph via a stable, unique substring from its original text.\nconst results = body.search(\"Acerca da quantidade de horas a serem alocadas\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  const hit = results.items[0];\n  const paragraph = hit.paragraphs.getFirst();\n  const paragraphRange = paragraph.getRange();\n\n  // Replacing the whole paragraph range collapses every run (and the\n  // spell-check proofErr markers around \"cargaHoraria\"/\"PRNCoordenador\"/\n  // \"coordenadorPRNTxtCoordenador\") into a single run that inherits the\n  // paragraph's run formatting, leaving only the final \"paragrafo8. \" text.\n  paragraphRange.insertText(\"paragrafo8. \", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Replace the paragraph that currently reads:\n#   \"Acerca da quantidade de horas a serem alocadas, por se tratar de um\n#    projeto cargaHoraria para PRNCoordenador coordenadorPRNTxtCoordenador,\n#    paragrafo8. \"\n# with just:\n#   \"paragrafo8. \"\n# keeping the paragraph's own formatting (indent / run size, etc.) intact.\n\n$d = $word.ActiveDocument\n\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($t -like \"*Acerca da quantidade de horas a serem alocadas*\") {\n        # Range that covers the paragraph text but not its trailing\n        # paragraph mark, so the replacement stays inside this paragraph.\n        $r = $d.Range($p.Range.Start, $p.Range.End - 1)\n        $r.Text = \"paragrafo8. \"\n        break\n    }\n}\n"}
